$d = $word.ActiveDocument

# Move to the very end of the document body (before the final paragraph mark)
$end = $d.Content.End - 1
$rng = $d.Range($end, $end)

# Insert a new paragraph after "Primeira vez escrevendo" with the new text
$rng.InsertParagraphAfter()
$rng.Collapse(0)
[void]$rng.MoveStart(1, 1)

$newPara = $d.Paragraphs.Last.Range
$newPara.Text = "Alteração que acho q é a certa"

# Word stamps the last edit position with a _GoBack bookmark;
# add it as a zero-length mark right after the text we just typed.
$bmRange = $d.Paragraphs.Last.Range
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
